$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, reusing the same formatting (style) as the
# other header cells, e.g. G1 ("sum").
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the data value for the new column in H2 (plain number, no special style)
$ws.Range("H2").Value = 0
